# Re-sort the worksheet tabs: "总计" (summary) moves in front of "2022-Q2"
# (the per-fund holdings detail), matching the commit's "resort sheetname"
# change. No cell data is modified - only tab order / which sheet is active.

$wb = $excel.ActiveWorkbook

# Move "总计" so it becomes the first tab (ahead of "2022-Q2").
$summarySheet = $wb.Worksheets.Item("总计")
$summarySheet.Move($wb.Worksheets.Item(1))

# Keep "2022-Q2" as the selected/active tab, same as before the resort.
# (Re-fetch by name instead of reusing a pre-move object reference, so the
# activation applies to the sheet in its new position.)
$wb.Worksheets.Item("2022-Q2").Activate()
